$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("G1").Value = "Situação"
$ws.Range("H1").Value = "Foi atendido"
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1:H1").PasteSpecial(-4122) | Out-Null

# Reorder existing rows: new row2 = old row4, new row3 = old row2, new row4 = old row3
$oldRow2 = @(1948252, "16/06/2021", 7565577, "AMBULATORIO SAUDE DO IDOSO", 3, 140035)
$oldRow3 = @(1948252, "16/06/2021", 7565577, "AMBULATORIO SAUDE DO IDOSO", 1, 140035)
$oldRow4 = @(1948257, "16/06/2021", 7565577, "AMBULATORIO SAUDE DO IDOSO", 1, 115853)

$newRows = @($oldRow4, $oldRow2, $oldRow3)

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $i + 2
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = "PACIENTE SEM SUSPEITA"
    $ws.Cells.Item($r, 8).Value = "NÃO"
}
